$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New participant row appended by the SmartScore Streamlit export (Arvand Zare)
$ws.Range("A15").Value = "Arvand Zare_20251202_124117"

# Grupo_Experimental was blank for this submission - force an empty text cell
# (not simply "no cell") to mirror the exported inline string, then drop the
# quote-prefix style picked up from the leading apostrophe trick.
$ws.Range("B15").Value = "'"
$ws.Range("B15").Style = "Normal"

$ws.Range("C15").Value = "Arvand Zare"
$ws.Range("D15").Value = 19
$ws.Range("E15").Value = "Male"
$ws.Range("F15").Value = "2025-12-02 12:41:18"
$ws.Range("G15").Value = "{`n  ""portion"": 0.6,`n  ""diet"": 1.0,`n  ""salt"": 0.8,`n  ""fat"": 0.6,`n  ""natural"": 0.8,`n  ""convenience"": 0.8,`n  ""price"": 0.4`n}"

$ws.Range("H15").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "0.523"
$ws.Range("I15").Style = "Normal"
$ws.Range("J15").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

$ws.Range("K15").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("L15").NumberFormat = "@"
$ws.Range("L15").Value = "0.457"
$ws.Range("L15").Style = "Normal"
$ws.Range("M15").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

$ws.Range("N15").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("O15").NumberFormat = "@"
$ws.Range("O15").Value = "0.427"
$ws.Range("O15").Style = "Normal"
$ws.Range("P15").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"

$ws.Range("Q15").Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "0.668"
$ws.Range("R15").Style = "Normal"
$ws.Range("S15").Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"

$ws.Range("T15").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("U15").NumberFormat = "@"
$ws.Range("U15").Value = "0.588"
$ws.Range("U15").Style = "Normal"
$ws.Range("V15").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

$ws.Range("W15").Value = "Annie’s Shells & White Cheddar"
$ws.Range("X15").NumberFormat = "@"
$ws.Range("X15").Value = "0.557"
$ws.Range("X15").Style = "Normal"
$ws.Range("Y15").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"

$ws.Range("Z15").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value = "0.778"
$ws.Range("AA15").Style = "Normal"
$ws.Range("AB15").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

$ws.Range("AC15").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AD15").NumberFormat = "@"
$ws.Range("AD15").Value = "0.591"
$ws.Range("AD15").Style = "Normal"
$ws.Range("AE15").Value = "Portátil, saludable, fácil, buena textura, sabor suave"

$ws.Range("AF15").Value = "Kitchens of India Variety Pack"
$ws.Range("AG15").NumberFormat = "@"
$ws.Range("AG15").Value = "0.559"
$ws.Range("AG15").Style = "Normal"
$ws.Range("AH15").Value = "Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad"
